# Fix email address in title slide: jdk72@zips.uakron.edu -> jdk72@uakron.edu
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldEmail = "jdk72@zips.uakron.edu"
$newEmail = "jdk72@uakron.edu"

foreach ($shp in $s.Shapes) {
    if ($shp.HasTextFrame) {
        $tf = $shp.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf($oldEmail)
            if ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $oldEmail.Length)
                $sub.Text = $newEmail
            }
        }
    }
}
